# Regenerate the handback-status report: refresh the "6c26c21b" /
# "ba65d440" pair's timestamps (and zh-cn priority) to reflect a newer
# report run. This fixture always keeps that pair's status columns in
# sync with each other (they were identical before the regen, and stay
# identical after it), so every cell sharing the old value is updated.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" column (G), rows 2 & 4
# (6c26c21b.md / ba65d440.md) 2016-09-04 00:18:32 -> 00:19:26
$wsOverview.Range("G2").Value = "2016-09-04 00:19:26"
$wsOverview.Range("G4").Value = "2016-09-04 00:19:26"

# zh-cn: rows 2 & 4 (6c26c21b / ba65d440)
#   Priority (E): ht -> mt
#   Correspond Handoff Datetime (H): 2016-09-04 00:18:27 -> 00:19:21
#   Correspond Handback DateTime (K): 2016-09-04 00:18:45 -> 00:19:42
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-04 00:19:21"
$wsZhCn.Range("H4").Value = "2016-09-04 00:19:21"
$wsZhCn.Range("K2").Value = "2016-09-04 00:19:42"
$wsZhCn.Range("K4").Value = "2016-09-04 00:19:42"

# de-de: rows 2 & 4 (6c26c21b / ba65d440)
#   Correspond Handback DateTime (K): 2016-09-04 00:18:52 -> 00:19:49
$wsDeDe.Range("K2").Value = "2016-09-04 00:19:49"
$wsDeDe.Range("K4").Value = "2016-09-04 00:19:49"
